$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'57.830.17"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").Value = "'3.133.32"
$ws.Range("E3").Value = "  +0.89%  "
$ws.Range("D5").Value = "'528.64"
$ws.Range("E5").Value = "  +0.96%  "
$ws.Range("D6").Value = "'138.52"
$ws.Range("E6").Value = "  -1.70%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "'3.134.97"
$ws.Range("E8").Value = "  +0.90%  "
$ws.Range("D9").Value = "'0.447"
$ws.Range("E9").Value = "  +3.02%  "
$ws.Range("E10").Value = "  -0.95%  "
$ws.Range("E11").Value = "  -1.12%  "
$ws.Range("E12").Value = "  +2.33%  "
$ws.Range("D13").Value = "'3.676.21"
$ws.Range("E13").Value = "  +0.99%  "
$ws.Range("E14").Value = "  +2.48%  "
$ws.Range("D15").Value = "'25.46"
$ws.Range("E15").Value = "  -3.01%  "
$ws.Range("D16").Value = "'0.0000164"
$ws.Range("E16").Value = "  +0.35%  "
$ws.Range("D17").Value = "'57.967.50"
$ws.Range("E17").Value = "  +0.46%  "
$ws.Range("D18").Value = "'3.158.63"
$ws.Range("E18").Value = "  +1.44%  "
$ws.Range("D19").Value = "'6.01"
$ws.Range("E19").Value = "  -1.86%  "
$ws.Range("D20").Value = "'12.78"
$ws.Range("E20").Value = "  -0.41%  "
$ws.Range("D21").Value = "'7.96"
$ws.Range("E21").Value = "  -1.50%  "
$ws.Range("D22").Value = "'353.94"
$ws.Range("E22").Value = "  +5.33%  "
$ws.Range("D23").Value = "'0.997"
$ws.Range("E23").Value = "  -0.28%  "
$ws.Range("D24").Value = "'68.85"
$ws.Range("E25").Value = "  -1.17%  "
$ws.Range("E26").Value = "  +0.24%  "
$ws.Range("D27").Value = "'0.998"
$ws.Range("E27").Value = "  -0.52%  "
$ws.Range("D28").Value = "'0.0₃0914"
$ws.Range("E28").Value = "  -1.07%  "
$ws.Range("E29").Value = "  +3.94%  "
$ws.Range("D31").Value = "'6.18"
$ws.Range("E31").Value = "  -5.67%  "
$ws.Range("D32").Value = "'1.88"
$ws.Range("E32").Value = "  +1.01%  "
$ws.Range("D33").Value = "'21.18"
$ws.Range("E33").Value = "  +0.90%  "
$ws.Range("D34").Value = "'1.18"
$ws.Range("E34").Value = "  -1.36%  "
$ws.Range("E35").Value = "  +6.77%  "
$ws.Range("D36").Value = "'158.87"
$ws.Range("E36").Value = "  +1.94%  "
$ws.Range("D37").Value = "'6.17"
$ws.Range("E37").Value = "  +1.09%  "
$ws.Range("D38").Value = "'26.72"
$ws.Range("E38").Value = "  -1.29%  "
$ws.Range("D39").Value = "'1.27"
$ws.Range("E39").Value = "  -1.90%  "
$ws.Range("D40").Value = "'0.0670"
$ws.Range("E40").Value = "  +1.04%  "
$ws.Range("E41").Value = "  +6.65%  "
$ws.Range("E42").Value = "  +5.88%  "
$ws.Range("E43").Value = "  +2.27%  "
$ws.Range("D44").Value = "'3.176.53"
$ws.Range("E44").Value = "  +0.88%  "
$ws.Range("D45").Value = "'0.0271"
$ws.Range("E45").Value = "  +4.94%  "
$ws.Range("D46").Value = "'36.54"
$ws.Range("E46").Value = "  -0.68%  "
$ws.Range("E47").Value = "  +0.03%  "
$ws.Range("D48").Value = "'2.312.33"
$ws.Range("E48").Value = "  +0.23%  "
$ws.Range("D49").Value = "'0.971"
$ws.Range("E49").Value = "  -0.64%  "
$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").Value = "'6.03"
$ws.Range("E50").Value = "  +0.30%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").Value = "'20.42"
$ws.Range("E51").Value = "  -1.56%  "
